$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the pFBA-related data rows into the transAnalysis sheet's data
# --- source table: add two new filename/metadata rows (8 and 9) just below
# --- the existing bop384 row, reusing the same column layout and style as
# --- the rows above them.

$ws.Range("A8").Value = "an_isoforms.fpkm_tracking"
$ws.Range("B8").Value = "RNAseq"
$ws.Range("C8").Value = "bop27?"
$ws.Range("D8").Value = "WT?"
$ws.Range("E8").Value = "Glucose?"
$ws.Range("F8").Value = "M9?"
$ws.Range("G8").Value = "Anaerobic"

$ws.Range("A9").Value = "arvsan_isoforms.fpkm_tracking"
$ws.Range("B9").Value = "RNAseq"
$ws.Range("C9").Value = "bop27?"
$ws.Range("D9").Value = "WT?"
$ws.Range("E9").Value = "Glucose?"
$ws.Range("F9").Value = "M9?"
$ws.Range("G9").Value = "Aerobic"

# Match the font/style used by the rest of the data table (rows 1-7): 12pt
# Verdana, the same as cell style index 1 applied to A1:H7.
$ws.Range("A8:G9").Font.Name = "Verdana"
$ws.Range("A8:G9").Font.Size = 12

# --- Move the current selection/cursor to G10, like the saved workbook.
$null = $ws.Range("G10").Select()

# --- Reposition the workbook window, matching the last-saved window state.
$excel.ActiveWindow.Left = 17840
$excel.ActiveWindow.Top = 3400
